# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and excel sheets
#
# 1) Trim now-irrelevant placeholder cells on "ODI Batting Extra" for rows
#    where the player did not bat at all (MATCH_CODE 4435 / 4436 / 4456):
#    only MATCH_CODE (A) and MAN_OF_MATCH (F) are kept.
# 2) Add a brand-new "ODI Bowling Extra" sheet (mirrors the structure of
#    "ODI Batting Extra" but for bowling: MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL)
#    and populate it with the scraped data for the last 20 ODI bowling
#    innings.

$wb = $excel.ActiveWorkbook

# --- 1) Clean up "ODI Batting Extra" -------------------------------------
$wsBattingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$wsBattingExtra.Range("B2:E2").ClearContents()
$wsBattingExtra.Range("B3:E3").ClearContents()
$wsBattingExtra.Range("B6:E6").ClearContents()

# --- 2) Add "ODI Bowling Extra" sheet -------------------------------------
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Header row, styled like the other sheets' header rows (bold, centered,
# top-aligned, thin border all round).
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL for the last
# 20 ODI bowling innings (same MATCH_CODEs as the tail of "ODI Bowling").
$rows = @(
    @{ Row=2;  A="4342"; B="0"; C=""       },
    @{ Row=3;  A="4345"; B="0"; C="30.00%" },
    @{ Row=4;  A="4350"; B="0"; C="10.00%" },
    @{ Row=5;  A="4353"; B="";  C=""       },
    @{ Row=6;  A="4436"; B="";  C=""       },
    @{ Row=7;  A="4457"; B="0"; C=""       },
    @{ Row=8;  A="4480"; B="0"; C="10.00%" },
    @{ Row=9;  A="4482"; B="0"; C=""       },
    @{ Row=10; A="4485"; B="0"; C="10.00%" },
    @{ Row=11; A="4609"; B="0"; C=""       },
    @{ Row=12; A="4613"; B="0"; C="20.00%" },
    @{ Row=13; A="4618"; B="3"; C="40.00%" },
    @{ Row=14; A="4687"; B="";  C=""       },
    @{ Row=15; A="4689"; B="0"; C=""       },
    @{ Row=16; A="4692"; B="";  C=""       },
    @{ Row=17; A="4695"; B="3"; C="20.00%" },
    @{ Row=18; A="4697"; B="0"; C="10.00%" },
    @{ Row=19; A="4725"; B="0"; C="10.00%" },
    @{ Row=20; A="4728"; B="0"; C=""       },
    @{ Row=21; A="4732"; B="0"; C="30.00%" }
)

foreach ($item in $rows) {
    $r = $item.Row

    # Every value is forced to text (leading apostrophe) so that
    # match-codes / "0" / percentages round-trip as strings rather than
    # being auto-coerced to numbers or percentage values, and so that
    # blank entries still materialise as an (empty) text cell instead of
    # being omitted entirely.
    $ws.Range("A$r").Value = "'" + $item.A

    if ($item.B -eq "") {
        $ws.Range("B$r").Value = "'"
    } else {
        $ws.Range("B$r").Value = "'" + $item.B
    }

    if ($item.C -eq "") {
        $ws.Range("C$r").Value = "'"
    } else {
        $ws.Range("C$r").Value = "'" + $item.C
    }
}
